# New crime data collected - weekly CompStat update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (Volume/Number and report week dates) ----
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# ---- Row 14 (Murder): D14/E14 go from numeric to "no data" markers ----
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("E23").Copy($ws.Range("E14"))

# ---- Row 15 (Rape) ----
$ws.Range("L15").Value = -33.333333333333

# ---- Row 16 (Robbery) ----
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 250
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 63.636363636363
$ws.Range("I16").Value = 192
$ws.Range("J16").Value = 240
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = 4.918032786885
$ws.Range("M16").Value = 32.413793103448
$ws.Range("N16").Value = -76.978417266187

# ---- Row 17 (Fel. Assault) ----
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 187
$ws.Range("J17").Value = 178
$ws.Range("K17").Value = 5.056179775280
$ws.Range("L17").Value = 10.650887573964
$ws.Range("M17").Value = 101.075268817204
$ws.Range("N17").Value = -37.458193979933

# ---- Row 18 (Burglary) ----
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = -46.153846153846
$ws.Range("I18").Value = 277
$ws.Range("J18").Value = 418
$ws.Range("K18").Value = -33.732057416267
$ws.Range("L18").Value = 18.376068376068
$ws.Range("M18").Value = 59.195402298850
$ws.Range("N18").Value = -64.758269720101

# ---- Row 19 (Gr. Larceny) ----
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 37
$ws.Range("E19").Value = -21.621621621621
$ws.Range("F19").Value = 97
$ws.Range("G19").Value = 111
$ws.Range("H19").Value = -12.612612612612
$ws.Range("I19").Value = 1283
$ws.Range("J19").Value = 1343
$ws.Range("K19").Value = -4.467609828741
$ws.Range("L19").Value = 51.654846335697
$ws.Range("M19").Value = 25.907752698724
$ws.Range("N19").Value = -48.618342010412

# ---- Row 20 (G.L.A.): D20/E20 go from numeric to "no data" markers ----
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E23").Copy($ws.Range("E20"))
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 0
$ws.Range("L20").Value = -18.181818181818
$ws.Range("M20").Value = 7.142857142857
$ws.Range("N20").Value = -94.023904382470

# ---- Row 21 (TOTAL) ----
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 164
$ws.Range("H21").Value = -14.024390243902
$ws.Range("I21").Value = 1993
$ws.Range("J21").Value = 2252
$ws.Range("K21").Value = -11.500888099467
$ws.Range("L21").Value = 32.955303535690
$ws.Range("M21").Value = 34.028244788164
$ws.Range("N21").Value = -61.539945966808

# ---- Row 22 (Transit) ----
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -62.5
$ws.Range("I22").Value = 42
$ws.Range("J22").Value = 45
$ws.Range("K22").Value = -6.666666666666
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = -14.285714285714

# ---- Row 24 (Petit Larceny) ----
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 3.125
$ws.Range("F24").Value = 135
$ws.Range("G24").Value = 136
$ws.Range("H24").Value = -0.735294117647
$ws.Range("I24").Value = 2004
$ws.Range("J24").Value = 1892
$ws.Range("K24").Value = 5.919661733615
$ws.Range("L24").Value = 58.544303797468
$ws.Range("M24").Value = 40.237928621413

# ---- Row 25 (Misd. Assault) ----
$ws.Range("C25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = 21.739130434782
$ws.Range("I25").Value = 410
$ws.Range("J25").Value = 403
$ws.Range("K25").Value = 1.736972704714
$ws.Range("L25").Value = 33.986928104575
$ws.Range("M25").Value = 64.658634538152

# ---- Row 26 (UCR Rape*): F26 goes from numeric to "no data" marker ----
$ws.Range("C14").Copy($ws.Range("F26"))
$ws.Range("H26").Value = -100
$ws.Range("I26").Value = 15
$ws.Range("K26").Value = -16.666666666666
$ws.Range("L26").Value = -6.25

# ---- Row 27 (Other Sex Crimes): D27/E27 go from numeric to "no data" markers ----
$ws.Range("C27").Value = 1
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E23").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 64
$ws.Range("K27").Value = -14.666666666666
$ws.Range("L27").Value = -5.882352941176

# ---- Row 30 (Hate Crimes): F30 goes from "no data" marker to numeric ----
$ws.Range("I30").Copy($ws.Range("F30"))
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 12
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -14.285714285714
